# The <id>p075v_N</id> markup for each of the 4 "tl" entries was
# originally split across three separate runs (one run for "<id>", a
# plainly-formatted run for "p075v_N", and another run for "</id>").
# Collapse each triple back into a single run carrying the "<id>"/"</id>"
# run's formatting (Courier New, color 7f6000, sz 18) by doing a
# find & replace of the already-identical visible text - Word merges the
# matched range into one run using the formatting of the first
# character of the match.
$d = $word.ActiveDocument
for ($i = 1; $i -le 4; $i++) {
    $old = "<id>p075v_$i</id>"
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $old, 2)
    Write-Host "Occurrence $i found:" $found
}
